$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Hoja1 (sheet1): re-point a few "KA5240_00" note cells in row 2 and drop
#    the leftover per-iteration notes in rows 3 and 4.
# ---------------------------------------------------------------------------
$ws1.Range("W2").Value = "KA5240_00_14993"
$ws1.Range("Z2").ClearContents()
$ws1.Range("Y2").Value = "Y luego por envios…."
$ws1.Range("W3").ClearContents()
$ws1.Range("W4").ClearContents()

# ---------------------------------------------------------------------------
# 2) Add the new "KA5240_00" sheet right after Hoja1 and make it the active
#    (selected) sheet/tab.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "KA5240_00"

$ws2.Range("A1").Value = "Envios"
$ws2.Range("B1").Value = "Public Score"

$envios = @(8000, 9000, 10000, 11000, 11500, 12000, 12500, 13000, 13500, 14000, 14500, 14993, 15000, 16000)
$scores = @(112, 124, 124, 125, 131, 135, 137, 135.7, 137.48, 134, 131.46, 129.78, 129, 131)

for ($i = 0; $i -lt $envios.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $envios[$i]
    $ws2.Cells.Item($r, 2).Value = $scores[$i]
}

# Highlight the chosen row (Envios = 13000) the way the author did, and leave
# it selected.
$ws2.Range("A9:B9").Interior.Color = 65535
$ws2.Range("A9:B9").Select()

# ---------------------------------------------------------------------------
# 3) Scatter chart of Public Score vs. Envios, embedded on the new sheet.
# ---------------------------------------------------------------------------
$chartObjs = $ws2.ChartObjects()
$co = $chartObjs.Add(80, 20, 430, 270)
$chart = $co.Chart
$chart.ChartType = 74  # xlXYScatterLines

$series = $chart.SeriesCollection().NewSeries()
$series.Name = "=KA5240_00!`$B`$1"
$series.XValues = $ws2.Range("A2:A15")
$series.Values = $ws2.Range("B2:B15")
$series.MarkerStyle = 8  # xlMarkerStyleCircle
$series.MarkerSize = 5

$chart.HasLegend = $false
$chart.HasTitle = $false

$xAxis = $chart.Axes(1)
$xAxis.MinimumScale = 8000

# ---------------------------------------------------------------------------
# 4) Leave Hoja1 selected cell where the author left it and switch focus to
#    the new sheet (KA5240_00 ends up as the active tab, like in the diff).
# ---------------------------------------------------------------------------
$ws1.Range("O7").Select()
$ws2.Activate()
